$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 1.7
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 2.3
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.9
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.75
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("X2").Value = 23
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 41
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 301
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 8
$ws.Range("AK2").Value = 13
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 26
$ws.Range("AP2").Value = 34
$ws.Range("AS2").Value = 251
$ws.Range("AT2").Value = 2.75
$ws.Range("AU2").Value = 8.5
$ws.Range("AX2").Value = 9
$ws.Range("AZ2").Value = 29
$ws.Range("BB2").Value = 151

$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 5.75
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 1.67
$ws.Range("K3").Value = 2.5
$ws.Range("L3").Value = 13
$ws.Range("W3").Value = 6
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 12
$ws.Range("AF3").Value = 151
$ws.Range("AH3").Value = 23
$ws.Range("AJ3").Value = 41
$ws.Range("AL3").Value = 126
$ws.Range("AW3").Value = 13
$ws.Range("AX3").Value = 67

$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.88

$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 2.9
$ws.Range("L6").Value = 3.5
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.95
$ws.Range("U6").Value = 1.73
$ws.Range("V6").Value = 2
$ws.Range("AK6").Value = 29
$ws.Range("AN6").Value = 4.5

$ws.Range("G7").Value = 3.3
$ws.Range("I7").Value = 2.25
$ws.Range("L7").Value = 3.1
$ws.Range("W7").Value = 8
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 41
$ws.Range("AA7").Value = 34
$ws.Range("AH7").Value = 6
$ws.Range("AI7").Value = 9.5
$ws.Range("AK7").Value = 21
$ws.Range("AL7").Value = 21
$ws.Range("AO7").Value = 21
$ws.Range("AS7").Value = 351
$ws.Range("AW7").Value = 4
$ws.Range("AX7").Value = 13

$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 2.75
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 4.75
$ws.Range("W8").Value = 6
$ws.Range("AA8").Value = 19
$ws.Range("AC8").Value = 8
$ws.Range("AE8").Value = 19
$ws.Range("AG8").Value = 501
$ws.Range("AH8").Value = 9
$ws.Range("AP8").Value = 26
$ws.Range("AU8").Value = 9
$ws.Range("AW8").Value = 6
$ws.Range("BA8").Value = 126

$ws.Range("G9").Value = 2.45
$ws.Range("I9").Value = 2.88
$ws.Range("Q9").Value = 2
$ws.Range("Y9").Value = 10
$ws.Range("AD9").Value = 6

$ws.Range("G10").Value = 1.6
$ws.Range("J10").Value = 2.2
$ws.Range("AE10").Value = 13
$ws.Range("AU10").Value = 7.5
$ws.Range("AX10").Value = 23

$ws.Range("G13").Value = 5.7
$ws.Range("H13").Value = 4.55
$ws.Range("I13").Value = 1.45
$ws.Range("J13").Value = 5.1
$ws.Range("K13").Value = 2.62
$ws.Range("L13").Value = 1.88
$ws.Range("N13").Value = 9.75
$ws.Range("R13").Value = 2.82
$ws.Range("S13").Value = 1.22
$ws.Range("T13").Value = 3.85
$ws.Range("U13").Value = 1.5
$ws.Range("V13").Value = 2.4
$ws.Range("W13").Value = 26
$ws.Range("X13").Value = 45
$ws.Range("Y13").Value = 18.5
$ws.Range("Z13").Value = 120
$ws.Range("AA13").Value = 45
$ws.Range("AB13").Value = 35
$ws.Range("AD13").Value = 10
$ws.Range("AE13").Value = 13.5
$ws.Range("AF13").Value = 40
$ws.Range("AG13").Value = 200
$ws.Range("AH13").Value = 11.5
$ws.Range("AI13").Value = 9.75
$ws.Range("AK13").Value = 12
$ws.Range("AL13").Value = 10.25
$ws.Range("AM13").Value = 17
$ws.Range("AN13").Value = 8
$ws.Range("AO13").Value = 28
$ws.Range("AP13").Value = 24
$ws.Range("AQ13").Value = 150
$ws.Range("AR13").Value = 120
$ws.Range("AS13").Value = 200
$ws.Range("AT13").Value = 3.85
$ws.Range("AU13").Value = 6.8
$ws.Range("AV13").Value = 40
$ws.Range("AW13").Value = 3.85
$ws.Range("AX13").Value = 6.7
$ws.Range("AZ13").Value = 17
$ws.Range("BA13").Value = 30
